# Applies the "Added a few more slots" edit to coywolf-cash review doc:
#   1. Remove the "Meta description: ..." paragraph that sits right after
#      the H1 title.
#   2. Insert a new bold paragraph re-stating the title ("Play Coywolf Cash
#      Free: Immersive American Wilderness Slot Game") right before the
#      final (italic) image-prompt paragraph.
#   3. Replace that final paragraph's italic text (the old image-generation
#      prompt) with the meta-description sentence that used to live at the
#      top of the document.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph (2nd paragraph) ------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# --- Step 2: insert a new bold paragraph just before the last paragraph ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount - 1)
$newRange = $newPara.Range

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr>' +
    '<w:t>Play Coywolf Cash Free: Immersive American Wilderness Slot Game</w:t>' +
    '</w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$newRange.InsertXML($titleXml)

# --- Step 3: swap the image-prompt text for the meta-description text ----
$oldText = "Please create a cartoon image for Coywolf Cash featuring a happy Maya warrior with glasses. The image should be fun and engaging, with bright colors and cartoon-style graphics. The Maya warrior should be smiling and holding a bag of money adorned with a dollar sign, with the Coywolf Cash slot machine in the background. The background of the image should feature the American wilderness, with rock formations, cacti and the endless road stretching out into the distance. The image should be eye-catching and encourage potential players to give Coywolf Cash a try."
$newText = "Read our game review for Coywolf Cash and play for free. Enjoy an immersive, American wilderness themed slot game with a high payout potential."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

Write-Output "Edit applied"
